# Split the single "Contributor" header column into two columns,
# mirroring the existing "Author First" / "Author Last" pattern, so
# contributor names can be captured (and later matched/recalled) the
# same way author names are.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing "Contributor" header (C1) and add the new "Contributor
# First" header in the previously-blank D1 cell.
$ws.Range("C1").Value = "Contributor Last"
$ws.Range("D1").Value = "Contributor First"

# Leave the cursor where the edit happened, matching the saved selection.
$ws.Range("D2").Select()
